# typewriter => release branch refresh
#
# Column G on Sheet1 lists "symbol" keyboard characters. Two new symbols
# (backtick and tilde) are inserted right under the "symbol" header, pushing
# every existing symbol down by two rows and growing the used range from
# A1:H31 to A1:H33. The row that used to carry the em dash (with its special
# Consolas-font style) moves along with its content to its new row, and the
# row that ends up holding the apostrophe picks up the quote-prefix style,
# same as the original file.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the Consolas-font style (currently on G22, the em dash) by
# stashing a copy of it onto G24 -- its new home -- before anything in
# column G gets overwritten with the shifted values.
$ws.Range("G22").Copy()
$ws.Range("G24").PasteSpecial(-4122)

# G22 itself becomes a plain "." in the new layout, so reset it back to the
# ordinary style shared by its neighbours (copy format from G21).
$ws.Range("G21").Copy()
$ws.Range("G22").PasteSpecial(-4122)

# New ordered list of symbols for G2:G33 (typewriter key-row order). Plain
# writes pick up the default style automatically.
#
# Two entries need special handling so Excel doesn't mis-parse them:
#  - a lone "'" would be swallowed as a pure quote-prefix marker (empty
#    text), so "''" is used instead, which stores a real apostrophe
#    character with the quote-prefix style switched on -- same result the
#    original author got by typing an apostrophe into a cell.
#  - a lone "=" would be parsed as an (empty) formula, so it is entered
#    as "'=" (quote-prefixed, forcing literal text); the quote-prefix style
#    that trick switches on is stripped right back off below.
$symbols = @(
    '`',
    '~',
    '!',
    '@',
    '#',
    '$',
    '%',
    '^',
    '&',
    '*',
    '(',
    ')',
    '-',
    "'=",
    '[',
    ']',
    ';',
    "''",
    '\',
    ',',
    '.',
    '/',
    '—',
    '+',
    '{',
    '}',
    ':',
    '"',
    '¦',
    '<',
    '>',
    '?'
)

$row = 2
foreach ($sym in $symbols) {
    $ws.Cells.Item($row, 7).Value = $sym
    $row = $row + 1
}

# Strip the incidental quote-prefix style picked up by the "'=" trick on
# row 15 (the "=" entry) -- copy plain formatting from its neighbour.
$ws.Range("G14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Selection moved from J19 to I8 in the saved view state.
$ws.Range("I8").Select()
